# Weekly update: a new "Rabanito" price record for
# "Vega Central Mapocho de Santiago" is inserted as the new row 139,
# pushing all the former rows 139-268 down by one (to 140-269).
# The sheet's used range grows from A1:R268 to A1:R269.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 139; Excel shifts rows 139..268
# down to 140..269 and extends the dimension automatically.
$ws.Rows.Item(139).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A139").Value = 9
$ws.Range("B139").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C139").Value = "Metropolitana"
$ws.Range("D139").Value = 44705
$ws.Range("E139").Value = 13
$ws.Range("F139").Value = 300000001
$ws.Range("G139").Value = "Rabanito"
$ws.Range("H139").Value = "Sin especificar"
$ws.Range("I139").Value = "Primera"
$ws.Range("J139").Value = 16000
$ws.Range("K139").Value = 2500
$ws.Range("L139").Value = 3000
$ws.Range("M139").Value = 2781
$ws.Range("N139").Value = "`$/cien unidades (volumen en unidades)"
$ws.Range("O139").Value = "Provincia de Chacabuco"
$ws.Range("P139").Value = 28
$ws.Range("Q139").Value = 100
$ws.Range("R139").Value = "Hortaliza"
